$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 18

$ws.Cells.Item($row - 1, 1).Copy() | Out-Null
$ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item($row, 1).Value = 42625.885057870371
$ws.Cells.Item($row, 2).Value = -30
$ws.Cells.Item($row, 3).Value = 63
$ws.Cells.Item($row, 4).Value = 35
$ws.Cells.Item($row, 5).Value = 71
$ws.Cells.Item($row, 6).Value = 28
$ws.Cells.Item($row, 7).Value = 11125
$ws.Cells.Item($row, 8).Value = 8723
$ws.Cells.Item($row, 9).Value = 461
$ws.Cells.Item($row, 10).Value = 132
$ws.Cells.Item($row, 11).Value = 74
$ws.Cells.Item($row, 12).Value = 10
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = "Named"
